$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking price strings
# (e.g. "1.00", "165.70") keep their exact textual representation instead
# of being coerced into numbers and losing formatting (trailing zeros, etc).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.385.62'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '2.578.81'
$ws.Range("E3").Value = '  -2.73%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '588.83'
$ws.Range("E5").Value = '  -3.17%  '
$ws.Range("D6").Value = '150.42'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").Value = '0.384'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '3.040.24'
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").Value = '63.204.06'
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("E16").Value = '  +5.17%  '
$ws.Range("D17").Value = '2.576.48'
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").Value = '12.18'
$ws.Range("E18").Value = '  +2.30%  '
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").Value = '  +3.26%  '
$ws.Range("D20").Value = '345.85'
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").Value = '67.06'
$ws.Range("E23").Value = '  +1.06%  '
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '9.14'
$ws.Range("E25").Value = '  -2.66%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '1.67'
$ws.Range("E26").Value = '  -3.46%  '
$ws.Range("D27").Value = '553.54'
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("D28").Value = '8.05'
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("D29").Value = '0.163'
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").Value = '2.03'
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").Value = '0.0₃0858'
$ws.Range("E32").Value = '  +0.65%  '
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").Value = '5.22'
$ws.Range("E34").Value = '  -1.15%  '
$ws.Range("D35").Value = '166.63'
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("D36").Value = '0.413'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").Value = '19.52'
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("E39").Value = '  -1.65%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").Value = '165.70'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").Value = '39.68'
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("D44").Value = '22.93'
$ws.Range("E44").Value = '  +3.60%  '
$ws.Range("D45").Value = '0.0585'
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("E46").Value = '  +5.27%  '
$ws.Range("D47").Value = '0.628'
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("D48").Value = '0.0251'
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").Value = '19.12'
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").Value = '0.0₆0235'
$ws.Range("E51").Value = '  +19.04%  '

# Reset column D style back to Normal (index 0) so we do not leave a
# lingering "@" text-format override on cells that did not need it.
$ws.Range("D2:D51").Style = "Normal"

